$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.621.25"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.226.77"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.50"
$ws.Range("E5").Value = "  -2.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "293.32"
$ws.Range("E6").Value = "  +10.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.66"
$ws.Range("E10").Value = "  -5.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0915"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.51"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.65"
$ws.Range("E13").Value = "  -6.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.05"
$ws.Range("E14").Value = "  +19.36%  "
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.96"
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.561.26"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.235.76"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.450.66"
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.21"
$ws.Range("E20").Value = "  +7.03%  "
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.54"
$ws.Range("E22").Value = "  +2.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.34"
$ws.Range("E23").Value = "  +14.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "235.32"
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("E26").Value = "  -4.77%  "
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.43"
$ws.Range("E28").Value = "  -8.28%  "
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.75"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.36"
$ws.Range("E31").Value = "  -8.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.12"
$ws.Range("E32").Value = "  -5.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.33"
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0882"
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.69"
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.02"
$ws.Range("E36").Value = "  +8.66%  "
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.75"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.29"
$ws.Range("E45").Value = "  -8.27%  "
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("E47").Value = "  -5.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.29"
$ws.Range("E48").Value = "  +3.24%  "
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.05"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.63"
$ws.Range("E51").Value = "  +3.95%  "
